$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C2:C5) from 2023-11-13 (45243) to 2023-11-14 (45244)
$ws.Range("C2").Value = 45244
$ws.Range("C3").Value = 45244
$ws.Range("C4").Value = 45244
$ws.Range("C5").Value = 45244
